# Update "想去人数" (F column) figures across sheets, matching
# the gh-pages data refresh described in the commit message
# ("Update gh-pages to output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1887
$ws.Range("F7").Value = 97
$ws.Range("F8").Value = 257
$ws.Range("F9").Value = 175
$ws.Range("F10").Value = 1117
$ws.Range("F18").Value = 112
$ws.Range("F19").Value = 1258
$ws.Range("F24").Value = 609
$ws.Range("F25").Value = 1040
$ws.Range("F26").Value = 59
$ws.Range("F27").Value = 1950
$ws.Range("F28").Value = 2417
$ws.Range("F33").Value = 562
$ws.Range("F34").Value = 735
$ws.Range("F35").Value = 793
$ws.Range("F38").Value = 733
$ws.Range("F39").Value = 218
$ws.Range("F40").Value = 551
$ws.Range("F41").Value = 641

# --- Sheet: 演出 (Performances) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F15").Value = 265
$ws.Range("F23").Value = 7

# --- Sheet: 全部类型 (All Types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1887
$ws.Range("F9").Value = 97
$ws.Range("F10").Value = 257
$ws.Range("F11").Value = 175
$ws.Range("F14").Value = 1117
$ws.Range("F22").Value = 112
$ws.Range("F23").Value = 1258
$ws.Range("F28").Value = 1040
$ws.Range("F29").Value = 2417
$ws.Range("F37").Value = 562
$ws.Range("F40").Value = 735
$ws.Range("F41").Value = 793
$ws.Range("F42").Value = 733
$ws.Range("F43").Value = 218
$ws.Range("F44").Value = 551
$ws.Range("F45").Value = 641
